# Trade #37 closed at 2026-02-17 12:39:11 - unknown UNKNOWN +0.000%
#
# Updates the Summary / Strategy Status sheets to reflect the new trade,
# and appends the closed-trade row to both "All Trades" and "MarketMaking".

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Summary sheet
# ---------------------------------------------------------------------
$summary = $wb.Worksheets.Item("Summary")
$summary.Range("B3").Value = 1200.72   # Current Capital
$summary.Range("B4").Value = 0.71      # Total P&L $
$summary.Range("B5").Value = 0.38      # Total P&L %
$summary.Range("B6").Value = 37        # Total Trades
$summary.Range("B8").Value = 14        # Losing Trades
$summary.Range("B9").Value = 37.84     # Win Rate %

# ---------------------------------------------------------------------
# Strategy Status sheet (MarketMaking row, row 4)
# ---------------------------------------------------------------------
$status = $wb.Worksheets.Item("Strategy Status")
$status.Range("C4").Value = 100.72     # Capital
$status.Range("D4").Value = 37         # Trades
$status.Range("E4").Value = 0.71       # P&L $
$status.Range("F4").Value = 0.72       # P&L %
$status.Range("G4").Value = 37.84      # Win Rate %

# ---------------------------------------------------------------------
# New closed-trade row, appended to both "All Trades" and "MarketMaking"
# ---------------------------------------------------------------------
$newRow = @(37, "2026-02-17", "12:39:05", "MarketMaking", "DOWN", 0.07000000000000001, 0.05, "CLOSED", -28.5714, -0.02, 100.72, 0, 0, 0.6, "Normal spread capture: 19600 bps", "early_exit", 0.13)

foreach ($sheetName in @("All Trades", "MarketMaking")) {
    $ws = $wb.Worksheets.Item($sheetName)
    $r = 38
    for ($i = 0; $i -lt $newRow.Length; $i++) {
        $col = $i + 1
        $cell = $ws.Cells.Item($r, $col)
        $val = $newRow[$i]

        if ($col -eq 2) {
            # "Date" column holds a plain text value ("2026-02-17") on every
            # other row (not a real Excel date). Entering that text straight
            # into .Value would make Excel auto-convert it to a date serial,
            # so force the cell to Text format first, then restore the
            # default style afterwards to leave no formatting footprint.
            $cell.NumberFormat = "@"
            $cell.Value = $val
            $cell.Style = "Normal"
        } else {
            $cell.Value = $val
        }
    }
}
